# Apply updated NATMI LR-pair data (H2-K1 -> Erbb2) per Dr Hou advice
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row2 = New-Object 'object[,]' 1,20
$row2[0,0] = "ECs"
$row2[0,1] = "H2-K1"
$row2[0,2] = "Erbb2"
$row2[0,3] = "ECs"
$row2[0,4] = 3
$row2[0,5] = 1
$row2[0,6] = 198.4680276666666
$row2[0,7] = 595.4040829999999
$row2[0,8] = 0.2835009389723355
$row2[0,9] = 0.2835009389723355
$row2[0,10] = 2
$row2[0,11] = 0.6666666666666666
$row2[0,12] = 1.720171333333333
$row2[0,13] = 5.160514
$row2[0,14] = 0.1961456356393658
$row2[0,15] = 0.1961456356393658
$row2[0,16] = 341.3990117754068
$row2[0,17] = 3072.591105978661
$row2[0,18] = 0.0556074718790858
$row2[0,19] = 0.0556074718790858
$ws.Range("A2:T2").Value = $row2

$row3 = New-Object 'object[,]' 1,20
$row3[0,0] = "ECs"
$row3[0,1] = "H2-K1"
$row3[0,2] = "Erbb2"
$row3[0,3] = "FAPs"
$row3[0,4] = 3
$row3[0,5] = 1
$row3[0,6] = 198.4680276666666
$row3[0,7] = 595.4040829999999
$row3[0,8] = 0.2835009389723355
$row3[0,9] = 0.2835009389723355
$row3[0,10] = 3
$row3[0,11] = 1
$row3[0,12] = 4.076388666666666
$row3[0,13] = 12.229166
$row3[0,14] = 0.464817562438416
$row3[0,15] = 0.464817562438416
$row3[0,16] = 809.0328186760863
$row3[0,17] = 7281.295368084777
$row3[0,18] = 0.1317762154021231
$row3[0,19] = 0.1317762154021231
$ws.Range("A3:T3").Value = $row3

$row4 = New-Object 'object[,]' 1,20
$row4[0,0] = "ECs"
$row4[0,1] = "H2-K1"
$row4[0,2] = "Erbb2"
$row4[0,3] = "sCs"
$row4[0,4] = 3
$row4[0,5] = 1
$row4[0,6] = 198.4680276666666
$row4[0,7] = 595.4040829999999
$row4[0,8] = 0.2835009389723355
$row4[0,9] = 0.2835009389723355
$row4[0,10] = 3
$row4[0,11] = 1
$row4[0,12] = 2.973308
$row4[0,13] = 8.919924
$row4[0,14] = 0.3390368019222182
$row4[0,15] = 0.3390368019222182
$row4[0,16] = 590.1065744055212
$row4[0,17] = 5310.959169649691
$row4[0,18] = 0.0961172516911266
$row4[0,19] = 0.09611725169112659
$ws.Range("A4:T4").Value = $row4

$row5 = New-Object 'object[,]' 1,20
$row5[0,0] = "FAPs"
$row5[0,1] = "H2-K1"
$row5[0,2] = "Erbb2"
$row5[0,3] = "ECs"
$row5[0,4] = 3
$row5[0,5] = 1
$row5[0,6] = 74.01845299999999
$row5[0,7] = 222.055359
$row5[0,8] = 0.1057313924740739
$row5[0,9] = 0.1057313924740739
$row5[0,10] = 2
$row5[0,11] = 0.6666666666666666
$row5[0,12] = 1.720171333333333
$row5[0,13] = 5.160514
$row5[0,14] = 0.1961456356393658
$row5[0,15] = 0.1961456356393658
$row5[0,16] = 127.3244209882807
$row5[0,17] = 1145.919788894526
$row5[0,18] = 0.02073875118386247
$row5[0,19] = 0.02073875118386247
$ws.Range("A5:T5").Value = $row5

$row6 = New-Object 'object[,]' 1,20
$row6[0,0] = "FAPs"
$row6[0,1] = "H2-K1"
$row6[0,2] = "Erbb2"
$row6[0,3] = "FAPs"
$row6[0,4] = 3
$row6[0,5] = 1
$row6[0,6] = 74.01845299999999
$row6[0,7] = 222.055359
$row6[0,8] = 0.1057313924740739
$row6[0,9] = 0.1057313924740739
$row6[0,10] = 3
$row6[0,11] = 1
$row6[0,12] = 4.076388666666666
$row6[0,13] = 12.229166
$row6[0,14] = 0.464817562438416
$row6[0,15] = 0.464817562438416
$row6[0,16] = 301.7279829333993
$row6[0,17] = 2715.551846400594
$row6[0,18] = 0.0491458081230185
$row6[0,19] = 0.0491458081230185
$ws.Range("A6:T6").Value = $row6

$row7 = New-Object 'object[,]' 1,20
$row7[0,0] = "FAPs"
$row7[0,1] = "H2-K1"
$row7[0,2] = "Erbb2"
$row7[0,3] = "sCs"
$row7[0,4] = 3
$row7[0,5] = 1
$row7[0,6] = 74.01845299999999
$row7[0,7] = 222.055359
$row7[0,8] = 0.1057313924740739
$row7[0,9] = 0.1057313924740739
$row7[0,10] = 3
$row7[0,11] = 1
$row7[0,12] = 2.973308
$row7[0,13] = 8.919924
$row7[0,14] = 0.3390368019222182
$row7[0,15] = 0.3390368019222182
$row7[0,16] = 220.079658452524
$row7[0,17] = 1980.716926072716
$row7[0,18] = 0.0358468331671929
$row7[0,19] = 0.03584683316719289
$ws.Range("A7:T7").Value = $row7

$row8 = New-Object 'object[,]' 1,20
$row8[0,0] = "M2"
$row8[0,1] = "H2-K1"
$row8[0,2] = "Erbb2"
$row8[0,3] = "ECs"
$row8[0,4] = 3
$row8[0,5] = 1
$row8[0,6] = 405.8333793333334
$row8[0,7] = 1217.500138
$row8[0,8] = 0.5797112283523728
$row8[0,9] = 0.5797112283523728
$row8[0,10] = 2
$row8[0,11] = 0.6666666666666666
$row8[0,12] = 1.720171333333333
$row8[0,13] = 5.160514
$row8[0,14] = 0.1961456356393658
$row8[0,15] = 0.1961456356393658
$row8[0,16] = 698.1029452389924
$row8[0,17] = 6282.926507150933
$row8[0,18] = 0.1137078273724537
$row8[0,19] = 0.1137078273724537
$ws.Range("A8:T8").Value = $row8

$row9 = New-Object 'object[,]' 1,20
$row9[0,0] = "M2"
$row9[0,1] = "H2-K1"
$row9[0,2] = "Erbb2"
$row9[0,3] = "FAPs"
$row9[0,4] = 3
$row9[0,5] = 1
$row9[0,6] = 405.8333793333334
$row9[0,7] = 1217.500138
$row9[0,8] = 0.5797112283523728
$row9[0,9] = 0.5797112283523728
$row9[0,10] = 3
$row9[0,11] = 1
$row9[0,12] = 4.076388666666666
$row9[0,13] = 12.229166
$row9[0,14] = 0.464817562438416
$row9[0,15] = 0.464817562438416
$row9[0,16] = 1654.334588069434
$row9[0,17] = 14889.01129262491
$row9[0,18] = 0.2694599600809299
$row9[0,19] = 0.2694599600809299
$ws.Range("A9:T9").Value = $row9

$row10 = New-Object 'object[,]' 1,20
$row10[0,0] = "M2"
$row10[0,1] = "H2-K1"
$row10[0,2] = "Erbb2"
$row10[0,3] = "sCs"
$row10[0,4] = 3
$row10[0,5] = 1
$row10[0,6] = 405.8333793333334
$row10[0,7] = 1217.500138
$row10[0,8] = 0.5797112283523728
$row10[0,9] = 0.5797112283523728
$row10[0,10] = 3
$row10[0,11] = 1
$row10[0,12] = 2.973308
$row10[0,13] = 8.919924
$row10[0,14] = 0.3390368019222182
$row10[0,15] = 0.3390368019222182
$row10[0,16] = 1206.667633438835
$row10[0,17] = 10860.00870094951
$row10[0,18] = 0.1965434408989892
$row10[0,19] = 0.1965434408989892
$ws.Range("A10:T10").Value = $row10

$row11 = New-Object 'object[,]' 1,20
$row11[0,0] = "sCs"
$row11[0,1] = "H2-K1"
$row11[0,2] = "Erbb2"
$row11[0,3] = "ECs"
$row11[0,4] = 3
$row11[0,5] = 1
$row11[0,6] = 21.741411
$row11[0,7] = 65.224233
$row11[0,8] = 0.03105644020121776
$row11[0,9] = 0.03105644020121776
$row11[0,10] = 2
$row11[0,11] = 0.6666666666666666
$row11[0,12] = 1.720171333333333
$row11[0,13] = 5.160514
$row11[0,14] = 0.1961456356393658
$row11[0,15] = 0.1961456356393658
$row11[0,16] = 37.398951948418
$row11[0,17] = 336.590567535762
$row11[0,18] = 0.006091585203963809
$row11[0,19] = 0.006091585203963809
$ws.Range("A11:T11").Value = $row11

$row12 = New-Object 'object[,]' 1,20
$row12[0,0] = "sCs"
$row12[0,1] = "H2-K1"
$row12[0,2] = "Erbb2"
$row12[0,3] = "FAPs"
$row12[0,4] = 3
$row12[0,5] = 1
$row12[0,6] = 21.741411
$row12[0,7] = 65.224233
$row12[0,8] = 0.03105644020121776
$row12[0,9] = 0.03105644020121776
$row12[0,10] = 3
$row12[0,11] = 1
$row12[0,12] = 4.076388666666666
$row12[0,13] = 12.229166
$row12[0,14] = 0.464817562438416
$row12[0,15] = 0.464817562438416
$row12[0,16] = 88.626441397742
$row12[0,17] = 797.6379725796779
$row12[0,18] = 0.01443557883234447
$row12[0,19] = 0.01443557883234447
$ws.Range("A12:T12").Value = $row12

$row13 = New-Object 'object[,]' 1,20
$row13[0,0] = "sCs"
$row13[0,1] = "H2-K1"
$row13[0,2] = "Erbb2"
$row13[0,3] = "sCs"
$row13[0,4] = 3
$row13[0,5] = 1
$row13[0,6] = 21.741411
$row13[0,7] = 65.224233
$row13[0,8] = 0.03105644020121776
$row13[0,9] = 0.03105644020121776
$row13[0,10] = 3
$row13[0,11] = 1
$row13[0,12] = 2.973308
$row13[0,13] = 8.919924
$row13[0,14] = 0.3390368019222182
$row13[0,15] = 0.3390368019222182
$row13[0,16] = 64.643911257588
$row13[0,17] = 581.795201318292
$row13[0,18] = 0.01052927616490948
$row13[0,19] = 0.01052927616490948
$ws.Range("A13:T13").Value = $row13
